$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update progress values from 80% to 90%
$ws.Range("C2").Value = 0.9
$ws.Range("C3").Value = 0.9

# Move the active selection from C5 to A8
$ws.Range("A8").Select()
